$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.Value = "'60.035.07"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.Value = "'  +2.05%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 4)
$c.Value = "'2.315.79"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.Value = "'  +0.32%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 5)
$c.Value = "'  +0.05%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 4)
$c.Value = "'542.26"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.Value = "'  +0.56%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.Value = "'130.58"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.Value = "'  -1.22%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.Value = "'  -0.01%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.576"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.Value = "'  -1.90%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.Value = "'2.315.37"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.Value = "'  +0.37%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.Value = "'  +0.22%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.Value = "'  +0.58%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.Value = "'  -0.07%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 4)
$c.Value = "'0.331"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.Value = "'  -0.84%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 4)
$c.Value = "'23.39"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.Value = "'  -1.49%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 4)
$c.Value = "'2.731.18"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.Value = "'  +0.29%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 4)
$c.Value = "'60.024.21"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.Value = "'  +2.23%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.Value = "'  -0.97%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.Value = "'2.317.59"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.Value = "'  -0.19%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 4)
$c.Value = "'10.50"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.Value = "'  -1.08%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 4)
$c.Value = "'4.09"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.Value = "'  -1.97%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 4)
$c.Value = "'312.76"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.Value = "'  -0.46%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.Value = "'  -0.93%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 4)
$c.Value = "'0.999"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.Value = "'  -0.13%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 2)
$c.Value = "'LEO"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 3)
$c.Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 4)
$c.Value = "'5.70"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 5)
$c.Value = "'  +0.56%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 2)
$c.Value = "'Litecoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 3)
$c.Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 4)
$c.Value = "'63.68"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.Value = "'  +1.42%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 2)
$c.Value = "'Kaspa"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 3)
$c.Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 4)
$c.Value = "'0.171"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.Value = "'  -0.79%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 2)
$c.Value = "'Binance-PegBSC-USD"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 3)
$c.Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 4)
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.Value = "'  +0.15%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 2)
$c.Value = "'InternetComputer(DFINITY)"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 3)
$c.Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 4)
$c.Value = "'7.75"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.Value = "'  -2.22%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 2)
$c.Value = "'Fetch.AI"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 3)
$c.Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 4)
$c.Value = "'1.34"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.Value = "'  +3.24%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 2)
$c.Value = "'SuiNetwork"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 3)
$c.Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 4)
$c.Value = "'1.19"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.Value = "'  +4.06%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 2)
$c.Value = "'Monero"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 3)
$c.Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 4)
$c.Value = "'171.16"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.Value = "'  +0.04%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 2)
$c.Value = "'PancakeSwap"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 3)
$c.Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 4)
$c.Value = "'1.71"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.Value = "'  -0.63%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 2)
$c.Value = "'PEPE"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 3)
$c.Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 4)
$c.Value = "'0.0₃0726"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.Value = "'  -1.15%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 2)
$c.Value = "'Aptos"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 3)
$c.Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 4)
$c.Value = "'5.86"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.Value = "'  -0.11%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 2)
$c.Value = "'ImmutableX"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 3)
$c.Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.Value = "'1.36"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.Value = "'  +4.90%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 2)
$c.Value = "'PolygonEcosystemToken"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 3)
$c.Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 4)
$c.Value = "'0.379"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.Value = "'  -1.72%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 2)
$c.Value = "'USDe"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 3)
$c.Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 4)
$c.Value = "'0.999"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 5)
$c.Value = "'  +0.01%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 2)
$c.Value = "'EthereumClassic"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 3)
$c.Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 4)
$c.Value = "'17.67"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.Value = "'  -1.46%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 2)
$c.Value = "'FirstDigitalUSD"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 3)
$c.Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 4)
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.Value = "'  +0.04%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 2)
$c.Value = "'NEARProtocol"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 3)
$c.Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 4)
$c.Value = "'4.01"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.Value = "'  -1.32%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 2)
$c.Value = "'Bittensor"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 3)
$c.Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.Value = "'317.95"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.Value = "'  +6.78%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 2)
$c.Value = "'OKB"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 3)
$c.Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 4)
$c.Value = "'37.86"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.Value = "'  -1.42%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 2)
$c.Value = "'Stacks"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 3)
$c.Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 4)
$c.Value = "'1.52"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.Value = "'  +0.10%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 2)
$c.Value = "'Aave"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 3)
$c.Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 4)
$c.Value = "'136.32"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.Value = "'  -3.79%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 2)
$c.Value = "'Filecoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 3)
$c.Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 4)
$c.Value = "'3.44"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.Value = "'  -0.22%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 2)
$c.Value = "'Stellar"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 3)
$c.Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 4)
$c.Value = "'0.0939"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.Value = "'  -2.40%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.Value = "'  +1.34%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 2)
$c.Value = "'InjectiveProtocol"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 3)
$c.Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 4)
$c.Value = "'18.78"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.Value = "'  +2.50%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 2)
$c.Value = "'Hedera"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 3)
$c.Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 4)
$c.Value = "'0.0490"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.Value = "'  -1.13%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 2)
$c.Value = "'BabyDogeCoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 3)
$c.Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 4)
$c.Value = "'0.0₆0227"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.Value = "'  +19.11%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 2)
$c.Value = "'VeChain"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 3)
$c.Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 4)
$c.Value = "'0.0211"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.Value = "'  +0.02%  "
$c.Style = "Normal"
